# "Partial fill out of QL form" — fill in the answers for the
# "JetBrains MPS" column (renamed from the template's "Language Workbench
# Name" placeholder) of the QL-features comparison table, plus a few
# "comment" cells in column C/E that explain specific answers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Header: name the language workbench being evaluated -------------
$ws.Range("B1").Value = "JetBrains MPS"

# --- Answers for each QL feature row ----------------------------------
# Cells that were fully blank (default style, just a thin border) lose
# that border once they are given a value, matching the template's
# "answered" look, so clear formatting before writing into them.
$answeredNoBorder = @("B2", "B4", "B5", "B6", "B9", "B15", "B16")
foreach ($ref in $answeredNoBorder) {
    $ws.Range($ref).Clear()
}

$ws.Range("B2").Value = "fully implemented"
$ws.Range("B4").Value = "fully implemented"
$ws.Range("B5").Value = "fully implemented"
$ws.Range("B6").Value = "not implemented"
$ws.Range("B9").Value = "fully implemented"
$ws.Range("B15").Value = "fully implemented"
$ws.Range("B16").Value = "fully implemented"

# Cells that already carried the bordered "table" style keep it — just
# set the value directly.
$ws.Range("B11").Value = "not implemented"
$ws.Range("C11").Value = "was asked not to implement as part of assignment"

$ws.Range("B18").Value = "not implemented"
$ws.Range("C18").Value = "was asked not to implement as part of assignment"

$ws.Range("B19").Value = "not implemented"
$ws.Range("C19").Value = "was asked not to implement as part of assignment"

$ws.Range("B20").Value = "not implemented"
$ws.Range("C20").Value = "was asked not to implement as part of assignment"

$ws.Range("B21").Value = "not implemented"
$ws.Range("C21").Value = "was asked not to implement as part of assignment"

$ws.Range("B22").Value = "not implemented"
$ws.Range("C22").Value = "was asked not to implement as part of assignment"

# New annotations in the little answer-key block (column E), echoing the
# legend entries already present in column E3:E7.
$ws.Range("E9").Value = "fully implemented"
$ws.Range("E10").Value = "partially implemented/limited support"
$ws.Range("E11").Value = "not implemented"

# B13 ("Highlighting" answer) goes back to a completely empty, unformatted
# cell (the legacy double-border placeholder style is removed).
$ws.Range("B13").Clear()

# --- Leave the cursor/selection where the author ended up editing ----
$ws.Range("B13").Select()
